$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1041.8572
$ws.Range("I2").Value = 215
$ws.Range("K2").Value = 215
$ws.Range("M2").Value = -102
# Row 19
$ws.Range("H19").Value = 55027.027
$ws.Range("I19").Value = 118642.7
$ws.Range("J19").Value = 3528.6191
$ws.Range("K19").Value = 118642.7
$ws.Range("L19").Value = 3528.6191
$ws.Range("M19").Value = -118467.7
$ws.Range("N19").Value = -3878.6191
# Row 33
$ws.Range("H33").Value = 33798.062
$ws.Range("I33").Value = 39751.383
$ws.Range("J33").Value = 8000.3335
$ws.Range("K33").Value = 39751.383
$ws.Range("L33").Value = 8000.3335
$ws.Range("M33").Value = -39522.383
$ws.Range("N33").Value = -8458.333500000001
# Row 40
$ws.Range("H40").Value = 4499.357
$ws.Range("I40").Value = 3262
$ws.Range("K40").Value = 3262
$ws.Range("M40").Value = -3087
# Row 82
$ws.Range("H82").Value = 10934.182
$ws.Range("I82").Value = 10036.25
$ws.Range("K82").Value = 30108.75
$ws.Range("M82").Value = -29702.75
# Row 85
$ws.Range("H85").Value = 10934.182
$ws.Range("I85").Value = 10036.25
$ws.Range("K85").Value = 30108.75
$ws.Range("M85").Value = -28704.75
# Row 86
$ws.Range("H86").Value = 7150.478
$ws.Range("I86").Value = 6305.625
$ws.Range("J86").Value = 7601.067
$ws.Range("K86").Value = 6305.625
$ws.Range("L86").Value = 7601.067
$ws.Range("M86").Value = -5182.625
$ws.Range("N86").Value = -9847.066999999999
# Row 89
$ws.Range("H89").Value = 7150.478
$ws.Range("I89").Value = 6305.625
$ws.Range("J89").Value = 7601.067
$ws.Range("K89").Value = 31528.125
$ws.Range("L89").Value = 38005.335
$ws.Range("M89").Value = -25912.125
$ws.Range("N89").Value = -49237.335
# Row 96
$ws.Range("H96").Value = 640.3182
$ws.Range("I96").Value = 415.2
$ws.Range("J96").Value = 827.9167
$ws.Range("K96").Value = 1245.6
$ws.Range("L96").Value = 2483.7501
$ws.Range("M96").Value = 127.4000000000001
$ws.Range("N96").Value = -5229.7501
# Row 135
$ws.Range("H135").Value = 2225.2173
$ws.Range("I135").Value = 1517.5
$ws.Range("J135").Value = 4773
$ws.Range("K135").Value = 13657.5
$ws.Range("L135").Value = 42957
$ws.Range("M135").Value = -11122.5
$ws.Range("N135").Value = -48027

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1315.8462
$ws.Range("I2").Value = 2205.8
$ws.Range("J2").Value = 759.625
$ws.Range("K2").Value = 2205.8
$ws.Range("L2").Value = 759.625
$ws.Range("M2").Value = -2092.8
$ws.Range("N2").Value = -985.625
# Row 97
$ws.Range("H97").Value = 211.23077
$ws.Range("I97").Value = 212.16667
$ws.Range("K97").Value = 212.16667
$ws.Range("M97").Value = 283.83333
# Row 116
$ws.Range("H116").Value = 1315.8462
$ws.Range("I116").Value = 2205.8
$ws.Range("J116").Value = 759.625
$ws.Range("K116").Value = 2205.8
$ws.Range("L116").Value = 759.625
$ws.Range("M116").Value = 88.19999999999982
$ws.Range("N116").Value = -5347.625
# Row 132
$ws.Range("H132").Value = 3616.6365
$ws.Range("I132").Value = 2148
$ws.Range("K132").Value = 6444
$ws.Range("M132").Value = -3914

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1315.8462
$ws.Range("I3").Value = 2205.8
$ws.Range("J3").Value = 759.625
$ws.Range("K3").Value = 2205.8
$ws.Range("L3").Value = 759.625
$ws.Range("M3").Value = -2091.8
$ws.Range("N3").Value = -987.625
# Row 21
$ws.Range("H21").Value = 68403.75
$ws.Range("J21").Value = 68403.75
$ws.Range("L21").Value = 68403.75
$ws.Range("N21").Value = -68875.75
# Row 86
$ws.Range("H86").Value = 7090.3
$ws.Range("I86").Value = 2474.5
$ws.Range("K86").Value = 2474.5
$ws.Range("M86").Value = -1351.5
# Row 89
$ws.Range("H89").Value = 7090.3
$ws.Range("I89").Value = 2474.5
$ws.Range("K89").Value = 12372.5
$ws.Range("M89").Value = -6756.5
# Row 134
$ws.Range("H134").Value = 3806
$ws.Range("I134").Value = 1778.2858
$ws.Range("K134").Value = 5334.857400000001
$ws.Range("M134").Value = -2799.857400000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 64999.5
$ws.Range("J18").Value = 64999.5
$ws.Range("L18").Value = 64999.5
$ws.Range("N18").Value = -65459.5
# Row 51
$ws.Range("H51").Value = 40555.375
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 40555.375
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 40555.375
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -42027.375
# Row 61
$ws.Range("H61").Value = 40555.375
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 40555.375
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 40555.375
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -41251.375
# Row 86
$ws.Range("H86").Value = 6001.8
$ws.Range("I86").Value = 4003
$ws.Range("K86").Value = 4003
$ws.Range("M86").Value = -2880
# Row 89
$ws.Range("H89").Value = 6001.8
$ws.Range("I89").Value = 4003
$ws.Range("K89").Value = 20015
$ws.Range("M89").Value = -14399
# Row 105
$ws.Range("H105").Value = 3264.72
$ws.Range("I105").Value = 2990.5
$ws.Range("K105").Value = 2990.5
$ws.Range("M105").Value = -1243.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 99
$ws.Range("H99").Value = 422.83334
$ws.Range("I99").Value = 447
$ws.Range("J99").Value = 418
$ws.Range("K99").Value = 1341
$ws.Range("L99").Value = 1254
$ws.Range("M99").Value = 905
$ws.Range("N99").Value = -5746
# Row 100
$ws.Range("H100").Value = 7599
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 6000
$ws.Range("M100").Value = -5189
# Row 133
$ws.Range("H133").Value = 14827.333
$ws.Range("J133").Value = 17777.666
$ws.Range("L133").Value = 53332.99800000001
$ws.Range("N133").Value = -63452.99800000001
# Row 141
$ws.Range("H141").Value = 3778.2727
$ws.Range("I141").Value = 3721.6667
$ws.Range("K141").Value = 11165.0001
$ws.Range("M141").Value = -5985.000100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5018.722
$ws.Range("I80").Value = 4481.143
$ws.Range("J80").Value = 5360.8184
$ws.Range("K80").Value = 4481.143
$ws.Range("L80").Value = 5360.8184
$ws.Range("M80").Value = -3483.143
$ws.Range("N80").Value = -7356.8184
# Row 83
$ws.Range("H83").Value = 5018.722
$ws.Range("I83").Value = 4481.143
$ws.Range("J83").Value = 5360.8184
$ws.Range("K83").Value = 22405.715
$ws.Range("L83").Value = 26804.092
$ws.Range("M83").Value = -17413.715
$ws.Range("N83").Value = -36788.092
# Row 94
$ws.Range("H94").Value = 179958.33
$ws.Range("J94").Value = 179958.33
$ws.Range("L94").Value = 179958.33
$ws.Range("N94").Value = -181310.33
# Row 97
$ws.Range("H97").Value = 439.14816
$ws.Range("I97").Value = 430.2381
$ws.Range("J97").Value = 470.33334
$ws.Range("K97").Value = 430.2381
$ws.Range("L97").Value = 470.33334
$ws.Range("M97").Value = 65.76190000000003
$ws.Range("N97").Value = -1462.33334
# Row 102
$ws.Range("H102").Value = 2474.2068
$ws.Range("I102").Value = 1470.88
$ws.Range("K102").Value = 1470.88
$ws.Range("M102").Value = 151.1199999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 20000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 20000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -20460
# Row 62
$ws.Range("H62").Value = 7713.5835
$ws.Range("J62").Value = 7438.7144
$ws.Range("L62").Value = 7438.7144
$ws.Range("N62").Value = -8686.714400000001
# Row 65
$ws.Range("H65").Value = 7713.5835
$ws.Range("J65").Value = 7438.7144
$ws.Range("L65").Value = 37193.572
$ws.Range("N65").Value = -43433.572
# Row 96
$ws.Range("H96").Value = 2112.7778
$ws.Range("I96").Value = 829.6667
$ws.Range("K96").Value = 829.6667
$ws.Range("M96").Value = 543.3333
